$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply new fill/highlight style (theme2, tint -0.249977111117893) to row 7 (A7:E7)
$ws.Range("A7:E7").Interior.ThemeColor = 2
$ws.Range("A7:E7").Interior.TintAndShade = -0.249977111117893

# Fill in row 9 with the new "Unbounded Knapsack" question data
$ws.Range("A9").Value = "CN"
$ws.Range("B9").Value = "CN"
$ws.Range("C9").Value = "Unbounded Knapsack"
$ws.Range("D9").Value = "Java"
$ws.Range("E9").Value = "DP(Recurrsion+Memonization+Tabulation+Space optimization)"

# Apply the same fill to row 9 and set row height
$ws.Range("A9:E9").Interior.ThemeColor = 2
$ws.Range("A9:E9").Interior.TintAndShade = -0.249977111117893
$ws.Rows.Item(9).RowHeight = 30

# Update selection to match C12:D12 with D12 active
$ws.Range("C12:D12").Select()
